$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before current row 432, restricted to the used columns
# (A:T) so we don't touch formatting for the whole 16384-column row.
$ws.Range("A432:T434").Insert(-4121)

# Copy formatting from the row that is now 435 (old row 432) down into the
# three newly inserted rows so styles (date format, borders, etc.) match.
$ws.Range("A435:T435").Copy()
$ws.Range("A432:T434").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 432
$ws.Range("A432").Value = 10
$ws.Range("B432").Value = "Vega Modelo de Temuco"
$ws.Range("C432").Value = "La Araucanía"
$ws.Range("D432").Value = 44656
$ws.Range("E432").Value = 9
$ws.Range("F432").Value = "Fruta"
$ws.Range("G432").Value = 100108
$ws.Range("H432").Value = "Tropicales y subtropicales"
$ws.Range("I432").Value = 100108005
$ws.Range("J432").Value = "Piña"
$ws.Range("K432").Value = "Caramelo"
$ws.Range("L432").Value = "Primera"
$ws.Range("M432").Value = 30
$ws.Range("N432").Value = 18000
$ws.Range("O432").Value = 18000
$ws.Range("P432").Value = 18000
$ws.Range("Q432").Value = "$/caja 12 unidades"
$ws.Range("R432").Value = "Ecuador"
$ws.Range("S432").Value = 1500
$ws.Range("T432").Value = 12

# New row 433
$ws.Range("A433").Value = 10
$ws.Range("B433").Value = "Vega Modelo de Temuco"
$ws.Range("C433").Value = "La Araucanía"
$ws.Range("D433").Value = 44656
$ws.Range("E433").Value = 9
$ws.Range("F433").Value = "Fruta"
$ws.Range("G433").Value = 100108
$ws.Range("H433").Value = "Tropicales y subtropicales"
$ws.Range("I433").Value = 100108005
$ws.Range("J433").Value = "Piña"
$ws.Range("K433").Value = "Caramelo"
$ws.Range("L433").Value = "Segunda"
$ws.Range("M433").Value = 30
$ws.Range("N433").Value = 18000
$ws.Range("O433").Value = 18000
$ws.Range("P433").Value = 18000
$ws.Range("Q433").Value = "$/caja 14 unidades"
$ws.Range("R433").Value = "Ecuador"
$ws.Range("S433").Value = 1286
$ws.Range("T433").Value = 14

# New row 434
$ws.Range("A434").Value = 10
$ws.Range("B434").Value = "Vega Modelo de Temuco"
$ws.Range("C434").Value = "La Araucanía"
$ws.Range("D434").Value = 44656
$ws.Range("E434").Value = 9
$ws.Range("F434").Value = "Fruta"
$ws.Range("G434").Value = 100108
$ws.Range("H434").Value = "Tropicales y subtropicales"
$ws.Range("I434").Value = 100108005
$ws.Range("J434").Value = "Piña"
$ws.Range("K434").Value = "Caramelo"
$ws.Range("L434").Value = "Tercera"
$ws.Range("M434").Value = 40
$ws.Range("N434").Value = 12000
$ws.Range("O434").Value = 12000
$ws.Range("P434").Value = 12000
$ws.Range("Q434").Value = "$/caja 16 unidades"
$ws.Range("R434").Value = "Ecuador"
$ws.Range("S434").Value = 750
$ws.Range("T434").Value = 16
